$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.87 = 6934.14 pesos`n✅ 6934.14 pesos = 1.86 = 944.89 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate table on "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 534.5
$ws2.Range("O10").Value = 3706.3
$ws2.Range("N12").Value = 3731.74
$ws2.Range("O12").Value = 508.508
